$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NutritionalData")
$ws.Range("J1").Value = '2 servings protein
(240	4	0	36	12	2	720)
3/4 cup mangos frozen
(90	0.5	0	1	20	2	0)
3/4 cups blueberries frozen
(52.5	0.75	0	0	14.25	3	0)
4 cups almond milk
(120	10	0	4	4	4	680)
8 corn tortillas Guerrero
(400	4	0	8	84	8	80)
1/2 cup mozzarella cheese
(160	10	7	12	2	0	380) 
bowl of beyond meat/3 bell peppers/2 zucchini
(290.33	24.67	4.33	7.67	8.00	2.00	119.67)
5 tbsp sourcream
(300	25	17.5	5	10	0	75)
nitro cold brew starbucks sweet cream
(70.00	5.00	3.50	1.00	4.00	0.00	15.00)
marguerite cauliflower pizza
(680	25	7	18	98	5	1040)
=240+90+52.5+120+400+160+290.33+300+70+680
=4+0.5+0.75+10+4+10+24.67+25+5+25
=0+0+0+0+0+7+4.33+17.5+3.5+7
=36+1+0+4+8+12+7.67+5+1+18
=12+20+14.25+4+84+2+8+10+4+98
=2+2+3+4+8+0+2+0+0+5
=720+0+0+680+80+380+119.67+75+15+1040'
